$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) DeploymentPlan sheet: swap the MAT_B@PLANT_001 block (rows 25-32, 8
#    rows) and the MAT_A@PLANT_001 block (rows 33-54, 22 rows) so the MAT_A
#    block now comes first, followed by the MAT_B block - i.e. rotate the
#    25-54 row range left by 8 rows.
# ---------------------------------------------------------------------------
$wsPlan = $wb.Worksheets.Item("DeploymentPlan")

$firstRow = 25
$lastRow = 54
$numCols = 13
$numRows = $lastRow - $firstRow + 1

# Snapshot every cell value in the block before overwriting anything.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    for ($c = 1; $c -le $numCols; $c++) {
        $snapshot[[string]$r + "_" + [string]$c] = $wsPlan.Cells.Item($r, $c).Value2
    }
}

# old row 33 -> new row 25 ... old row 54 -> new row 46
# old row 25 -> new row 47 ... old row 32 -> new row 54
$rotateBy = 8
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $oldRow = $r + $rotateBy
    if ($oldRow -gt $lastRow) {
        $oldRow = $oldRow - $numRows
    }
    for ($c = 1; $c -le $numCols; $c++) {
        $wsPlan.Cells.Item($r, $c).Value = $snapshot[[string]$oldRow + "_" + [string]$c]
    }
}

# ---------------------------------------------------------------------------
# 2) StockOnHandLog sheet: reorder rows 3-7 (row 2 is unchanged).
#    new row 3 <- old row 5
#    new row 4 <- old row 7
#    new row 5 <- old row 6
#    new row 6 <- old row 4
#    new row 7 <- old row 3
# ---------------------------------------------------------------------------
$wsStock = $wb.Worksheets.Item("StockOnHandLog")

$stockNumCols = 10
$stockOrder = @{ 3 = 5; 4 = 7; 5 = 6; 6 = 4; 7 = 3 }

$stockSnapshot = @{}
foreach ($r in 3..7) {
    for ($c = 1; $c -le $stockNumCols; $c++) {
        $stockSnapshot[[string]$r + "_" + [string]$c] = $wsStock.Cells.Item($r, $c).Value2
    }
}

foreach ($newRow in 3..7) {
    $oldRow = $stockOrder[$newRow]
    for ($c = 1; $c -le $stockNumCols; $c++) {
        $wsStock.Cells.Item($newRow, $c).Value = $stockSnapshot[[string]$oldRow + "_" + [string]$c]
    }
}
